$d = $word.ActiveDocument

# Old wording (as currently split across two runs, joined by a manual line break)
$oldPart1 = "## Warning in (function (z, notch = FALSE, width = NULL, varwidth = FALSE, : some notches went outside"
$oldPart2 = "## hinges ('box'): maybe set notch=FALSE"

# New wording: " hinges ('box'):" moves up to the end of the first line
$newPart1 = "## Warning in (function (z, notch = FALSE, width = NULL, varwidth = FALSE, : some notches went outside hinges ('box'):"
$newPart2 = "## maybe set notch=FALSE"

# Locate the paragraph that holds this warning text.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$oldPart1*") {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    # Replacing the whole-paragraph range's .Text with the new first-line text only
    # rewrites the first run in place (keeping its rPr/rStyle and avoiding the
    # curly-quote autocorrect that Find/Replace would otherwise trigger), and
    # leaves the manual line break plus the second run completely untouched.
    $paraRange = $targetPara.Range.Duplicate
    $paraRange.Text = $newPart1

    # Now locate the (still original) second run's text so we can replace it too.
    $searchRange = $targetPara.Range.Duplicate
    $found = $searchRange.Find.Execute($oldPart2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $runRange = $searchRange.Duplicate
        $runRange.Text = $newPart2
    }
}
